# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The account-statement table (rows 16-22, columns E "Periodo Mora" and
# F "Valor Mora") is reordered: the period list is reversed so that the
# most recent period (2501) now appears first and the oldest (2407) now
# appears last, with each period's "Valor Mora" value travelling together
# with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of periods (reversed) for rows 16..22
$periods = @("2501", "2412", "2411", "2410", "2409", "2408", "2407")

# Valor Mora that travels together with its period (old row 22 had 50266,
# everything else was 52000); after reversing, row 16 (period 2501) now
# carries 50266 and row 22 (period 2407) now carries 52000.
$valores = @(50266, 52000, 52000, 52000, 52000, 52000, 52000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
